$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRODUCTO")

# New H-column values (3rd parameter passed to PG_CI_PUNTO_VENTA) for the
# "SK_Punto Venta Completo" rebalance. Row -> new value.
$changes = [ordered]@{
    9  = 2
    10 = 3
    11 = 1
    13 = 2
    15 = 3
    18 = 2
    20 = 3
    21 = 1
    23 = 2
    25 = 3
    28 = 2
    30 = 3
    31 = 1
    33 = 3
    34 = 2
    35 = 3
    38 = 2
    40 = 3
    41 = 1
    43 = 2
    45 = 3
    48 = 2
    49 = 2
}

foreach ($row in $changes.Keys) {
    $cell = $ws.Range("H$row")
    $cell.Value = $changes[$row]
    # Re-apply the column's normal number format (quote-prefixed "General")
    # so the cell keeps its original style instead of reverting to the
    # workbook default after the value write.
    $ws.Range("H3").Copy()
    $cell.PasteSpecial(-4122)
}

# Row 24 keeps its value (3) but is now flagged/underlined for review.
$ws.Range("H24").Font.Underline = $true

# The user scrolled the frozen pane down and selected the last data row.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("H51").Select()
